$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the current row 114, pushing the
# existing rows 115-134 down to 117-136 (mirrors the weekly data refresh
# described in the commit message: two new observations were added).
$ws.Rows("115:116").Insert()

# New row 115: Apio / Americana (o) / Primera, week of 44476
$ws.Range("A115").Value = 9
$ws.Range("B115").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C115").Value = "Metropolitana"
$ws.Range("D115").Value = 44476
$ws.Range("E115").Value = 13
$ws.Range("F115").Value = 100112017
$ws.Range("G115").Value = "Apio"
$ws.Range("H115").Value = "Americana (o)"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 61
$ws.Range("K115").Value = 8000
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = 8492
$ws.Range("N115").Value = '$/docena de matas'
$ws.Range("O115").Value = "Región de Coquimbo"
$ws.Range("P115").Value = 1415
$ws.Range("Q115").Value = 6
$ws.Range("R115").Value = "Hortaliza"

# New row 116: Apio / Americana (o) / Segunda, week of 44476
$ws.Range("A116").Value = 9
$ws.Range("B116").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C116").Value = "Metropolitana"
$ws.Range("D116").Value = 44476
$ws.Range("E116").Value = 13
$ws.Range("F116").Value = 100112017
$ws.Range("G116").Value = "Apio"
$ws.Range("H116").Value = "Americana (o)"
$ws.Range("I116").Value = "Segunda"
$ws.Range("J116").Value = 43
$ws.Range("K116").Value = 6000
$ws.Range("L116").Value = 7000
$ws.Range("M116").Value = 6512
$ws.Range("N116").Value = '$/docena de matas'
$ws.Range("O116").Value = "Región de Coquimbo"
$ws.Range("P116").Value = 1085
$ws.Range("Q116").Value = 6
$ws.Range("R116").Value = "Hortaliza"

# Make sure column D keeps the date format used by every other row in
# this column.
$ws.Range("D115:D116").NumberFormat = "YYYY-MM-DD HH:MM:SS"
